# Rows 50-55 of the "Artfynd" sheet need to have their data cyclically
# shifted down by one row, with the content that was on row 55 wrapping
# around to become the new row 50 (new_row[51..55] = old_row[50..54],
# new_row[50] = old_row[55]). The row numbers / record IDs stay put; only
# the field values move between rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcRange = "A50:AY55"

$rng = $ws.Range($srcRange)
$vals = $rng.Value2

$rows = $vals.GetLength(0)
$cols = $vals.GetLength(1)

# Build the rotated block in a new (0-based) array: row 1 gets what used to
# be in the last row, and every other row gets what used to be one row above.
$newVals = New-Object 'object[,]' $rows,$cols

for ($c = 1; $c -le $cols; $c++) {
    $newVals[0, $c - 1] = $vals[$rows, $c]
    for ($r = 2; $r -le $rows; $r++) {
        $newVals[$r - 1, $c - 1] = $vals[$r - 1, $c]
    }
}

# The Startdatum/Slutdatum columns (Y and AA) hold dates stored as plain
# text (e.g. "2023-09-06"). A bulk Value2 write auto-detects such strings
# and silently turns them into real date serials, which would corrupt the
# data. Temporarily force those columns to Text format so the values are
# written back verbatim as strings.
$ws.Range("Y50:Y55").NumberFormat = "@"
$ws.Range("AA50:AA55").NumberFormat = "@"

$ws.Range($srcRange).Value2 = $newVals

# Restore the default/implicit styling for the whole block (none of these
# cells had an explicit style originally) so we don't leave a stray "Text"
# number format behind on the date columns.
$ws.Range($srcRange).Style = "Normal"
